# Auto-generated script to apply cryptos.xlsx price/volume update
# Commit: Updated cryptos list on Sat Jun 24 04:18:51 UTC 2023 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws 'D2' '30.713.64'
Set-TextCell $ws 'E2' '  +2.50%  '
Set-TextCell $ws 'D3' '1.896.25'
Set-TextCell $ws 'E3' '  +0.96%  '
Set-TextCell $ws 'D4' '1.001'
Set-TextCell $ws 'E4' '  +0.30%  '
Set-TextCell $ws 'D5' '245.97'
Set-TextCell $ws 'E5' '  +1.80%  '
Set-TextCell $ws 'D6' '1.000'
Set-TextCell $ws 'E6' '  +0.29%  '
Set-TextCell $ws 'D7' '0.4922'
Set-TextCell $ws 'E7' '  -1.46%  '
Set-TextCell $ws 'D8' '0.2951'
Set-TextCell $ws 'E8' '  +0.94%  '
Set-TextCell $ws 'D9' '0.06824'
Set-TextCell $ws 'E9' '  +3.19%  '
Set-TextCell $ws 'D10' '1.889.73'
Set-TextCell $ws 'E10' '  +0.63%  '
Set-TextCell $ws 'D11' '17.34'
Set-TextCell $ws 'E11' '  +3.64%  '
Set-TextCell $ws 'B12' 'TRON'
Set-TextCell $ws 'C12' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell $ws 'D12' '0.07267'
Set-TextCell $ws 'E12' '  +0.08%  '
Set-TextCell $ws 'B13' 'Litecoin'
Set-TextCell $ws 'C13' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell $ws 'D13' '92.24'
Set-TextCell $ws 'E13' '  +7.12%  '
Set-TextCell $ws 'B14' 'Polkadot'
Set-TextCell $ws 'C14' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell $ws 'D14' '5.121'
Set-TextCell $ws 'E14' '  +5.41%  '
Set-TextCell $ws 'B15' 'Polygon'
Set-TextCell $ws 'C15' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell $ws 'D15' '0.6845'
Set-TextCell $ws 'E15' '  +2.68%  '
Set-TextCell $ws 'D16' '30.684.26'
Set-TextCell $ws 'E16' '  +2.50%  '
Set-TextCell $ws 'D17' '0.000007992'
Set-TextCell $ws 'E17' '  +1.26%  '
Set-TextCell $ws 'D18' '13.38'
Set-TextCell $ws 'E18' '  +4.86%  '
Set-TextCell $ws 'D19' '1.0000'
Set-TextCell $ws 'E19' '  +0.26%  '
Set-TextCell $ws 'D20' '2.141.06'
Set-TextCell $ws 'E20' '  +1.08%  '
Set-TextCell $ws 'D21' '0.9997'
Set-TextCell $ws 'E21' '  +0.24%  '
Set-TextCell $ws 'D22' '4.870'
Set-TextCell $ws 'E22' '  +2.45%  '
Set-TextCell $ws 'D23' '189.82'
Set-TextCell $ws 'E23' '  +36.23%  '
Set-TextCell $ws 'D24' '6.089'
Set-TextCell $ws 'E24' '  +8.02%  '
Set-TextCell $ws 'D25' '9.361'
Set-TextCell $ws 'E25' '  +3.36%  '
Set-TextCell $ws 'D26' '155.14'
Set-TextCell $ws 'E26' '  +3.97%  '
Set-TextCell $ws 'D27' '19.30'
Set-TextCell $ws 'E27' '  +14.00%  '
Set-TextCell $ws 'D28' '1.934'
Set-TextCell $ws 'E28' '  +1.29%  '
Set-TextCell $ws 'D29' '1.397'
Set-TextCell $ws 'E29' '  +0.70%  '
Set-TextCell $ws 'D30' '4.376'
Set-TextCell $ws 'E30' '  +4.68%  '
Set-TextCell $ws 'D31' '0.09023'
Set-TextCell $ws 'E31' '  +2.64%  '
Set-TextCell $ws 'D32' '4.064'
Set-TextCell $ws 'E32' '  +2.99%  '
Set-TextCell $ws 'D33' '0.05187'
Set-TextCell $ws 'E33' '  +2.60%  '
Set-TextCell $ws 'D34' '0.7494'
Set-TextCell $ws 'E34' '  +4.64%  '
Set-TextCell $ws 'D35' '1.129'
Set-TextCell $ws 'E35' '  +2.01%  '
Set-TextCell $ws 'D36' '2.710'
Set-TextCell $ws 'E36' '  +1.74%  '
Set-TextCell $ws 'D37' '0.01876'
Set-TextCell $ws 'E37' '  +7.09%  '
Set-TextCell $ws 'D38' '2.672'
Set-TextCell $ws 'E38' '  -0.81%  '
Set-TextCell $ws 'D39' '2.170'
Set-TextCell $ws 'E39' '  -0.50%  '
Set-TextCell $ws 'D40' '0.9378'
Set-TextCell $ws 'D41' '0.4445'
Set-TextCell $ws 'E41' '  +4.41%  '
Set-TextCell $ws 'D42' '106.22'
Set-TextCell $ws 'E42' '  +4.46%  '
Set-TextCell $ws 'D43' '5.813'
Set-TextCell $ws 'E43' '  +0.37%  '
Set-TextCell $ws 'D44' '1.001'
Set-TextCell $ws 'E44' '  +0.37%  '
Set-TextCell $ws 'D45' '7.727'
Set-TextCell $ws 'E45' '  +3.66%  '
Set-TextCell $ws 'D46' '0.1346'
Set-TextCell $ws 'E46' '  +7.18%  '
Set-TextCell $ws 'D47' '0.05858'
Set-TextCell $ws 'E47' '  +3.59%  '
Set-TextCell $ws 'D48' '8.786'
Set-TextCell $ws 'E48' '  +7.46%  '
Set-TextCell $ws 'D49' '0.3966'
Set-TextCell $ws 'E49' '  +5.54%  '
Set-TextCell $ws 'B50' 'NEARProtocol'
Set-TextCell $ws 'C50' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws 'D50' '1.430'
Set-TextCell $ws 'E50' '  +7.41%  '
Set-TextCell $ws 'B51' 'Elrond'
Set-TextCell $ws 'C51' 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextCell $ws 'D51' '33.63'
Set-TextCell $ws 'E51' '  +3.91%  '
